$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is a duplicate of row 3 (Ernie's record) except column A gets a new,
# deliberately very long "alias" name -- added to see how a long driver name
# affects the layout/format (per the commit message).
$ws.Range("A4").Value = "ErnieAliasBertiBertLongNameTest"
$ws.Range("B4").Value = $ws.Range("B3").Value2
$ws.Range("C4").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("D3").Value2
$ws.Range("E4").Value = $ws.Range("E3").Value2
$ws.Range("F4").Value = $ws.Range("F3").Value2
$ws.Range("G4").Value = $ws.Range("G3").Value2

# Nudge the font on the long-name cell (and on A8, further down the sheet)
# so a dedicated style/font entry gets recorded for them.
$ws.Range("A4").Font.ColorIndex = 1
$ws.Range("A8").Font.ColorIndex = 1

# Leave the selection where the author ended up after poking at the new
# layout.
$ws.Range("A8").Select()
